$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 29; $r -le 36; $r++) {
    $ws.Range("F$r").Formula = "=ROUND(E$r*`$D$r,0)"
    $ws.Range("H$r").Formula = "=ROUND(G$r*`$D$r,0)"
    $ws.Range("J$r").Formula = "=ROUND(I$r*`$D$r,0)"
    $ws.Range("L$r").Formula = "=ROUND(K$r*`$D$r,0)"
    $ws.Range("N$r").Formula = "=ROUND(M$r*`$D$r,0)"
    $ws.Range("P$r").Formula = "=ROUND(O$r*`$D$r,0)"
}

$ws.Range("R33").Select()
